$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# New row 38: "cpi" label in column A (same style/format as the "inflation" row above it)
$ws.Cells.Item(38, 1).Value = "cpi"

# Columns E:R (1920-1980) hold a flat 116.4 placeholder value, shown in red
# to flag that it is not real historical data but a held-constant filler.
$constCols = 5..18
foreach ($col in $constCols) {
    $cell = $ws.Cells.Item(38, $col)
    $cell.Value = 116.4
    $cell.Font.Color = 255
}

# Columns S:BA (1981-2015) hold the actual CPI index series (base year 2010 = 100).
$cpiValues = @(28.8,30.8,32.9,34.6,36.1,36.8,38.9,43.4,46.1,49.9,53.1,55.7,58.6,59.1,60.4,62.2,64.3,68,69.3,70.1,73.1,73.6,76.7,78.9,80.8,82.8,86.7,93.4,96.5,100,107.3,110.9,114.2,114.8,116.4)
$startCol = 19
for ($i = 0; $i -lt $cpiValues.Length; $i++) {
    $ws.Cells.Item(38, $startCol + $i).Value = $cpiValues[$i]
}

# Keep the view scrolled near the newly added row, matching the author's
# navigation at the time of the edit, and select the new row label.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A38").Select()
